$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows 3, 4, 5 and 7 have their Fecha/Variedad/Calidad/Precio
# mínimo/Precio máximo/Precio promedio ponderado/Precio $/Kg values
# cyclically rotated: old row3 -> row4, old row4 -> row7, old row7 -> row5,
# old row5 -> row3. Capture the "before" values first, then write them
# back into their new rows.

$row3 = @{
    D = $ws.Range("D3").Value()
    K = $ws.Range("K3").Value()
    L = $ws.Range("L3").Value()
    N = $ws.Range("N3").Value()
    O = $ws.Range("O3").Value()
    P = $ws.Range("P3").Value()
    S = $ws.Range("S3").Value()
}

$row4 = @{
    D = $ws.Range("D4").Value()
    K = $ws.Range("K4").Value()
    L = $ws.Range("L4").Value()
    N = $ws.Range("N4").Value()
    O = $ws.Range("O4").Value()
    P = $ws.Range("P4").Value()
    S = $ws.Range("S4").Value()
}

$row5 = @{
    D = $ws.Range("D5").Value()
    K = $ws.Range("K5").Value()
    L = $ws.Range("L5").Value()
    N = $ws.Range("N5").Value()
    O = $ws.Range("O5").Value()
    P = $ws.Range("P5").Value()
    S = $ws.Range("S5").Value()
}

$row7 = @{
    D = $ws.Range("D7").Value()
    K = $ws.Range("K7").Value()
    L = $ws.Range("L7").Value()
    N = $ws.Range("N7").Value()
    O = $ws.Range("O7").Value()
    P = $ws.Range("P7").Value()
    S = $ws.Range("S7").Value()
}

# New row 3 = old row 5
$ws.Range("D3").Value = $row5.D
$ws.Range("K3").Value = $row5.K
$ws.Range("L3").Value = $row5.L
$ws.Range("N3").Value = $row5.N
$ws.Range("O3").Value = $row5.O
$ws.Range("P3").Value = $row5.P
$ws.Range("S3").Value = $row5.S

# New row 4 = old row 3
$ws.Range("D4").Value = $row3.D
$ws.Range("K4").Value = $row3.K
$ws.Range("L4").Value = $row3.L
$ws.Range("N4").Value = $row3.N
$ws.Range("O4").Value = $row3.O
$ws.Range("P4").Value = $row3.P
$ws.Range("S4").Value = $row3.S

# New row 5 = old row 7
$ws.Range("D5").Value = $row7.D
$ws.Range("K5").Value = $row7.K
$ws.Range("L5").Value = $row7.L
$ws.Range("N5").Value = $row7.N
$ws.Range("O5").Value = $row7.O
$ws.Range("P5").Value = $row7.P
$ws.Range("S5").Value = $row7.S

# New row 7 = old row 4
$ws.Range("D7").Value = $row4.D
$ws.Range("K7").Value = $row4.K
$ws.Range("L7").Value = $row4.L
$ws.Range("N7").Value = $row4.N
$ws.Range("O7").Value = $row4.O
$ws.Range("P7").Value = $row4.P
$ws.Range("S7").Value = $row4.S
